$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (2-485).
# This refresh run bumps that "last changed" stamp from 2023-09-17 (45186)
# to 2023-09-19 (45188) for all of them.
$ws.Range("C2:C485").Value = 45188
